$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the user/email rows: "user03"/"tu03@maildrop.cc" -> "user09"/"tu09@maildrop.cc"
# and "user04"/"tu04@maildrop.cc" -> "user08"/"tu08@maildrop.cc"
$ws.Range("B2").Value = "user09"
$ws.Range("C2").Value = "tu09@maildrop.cc"
$ws.Range("B3").Value = "user08"
$ws.Range("C3").Value = "tu08@maildrop.cc"

# Update selection to C3
$ws.Range("C3").Select()
